# Update the day/time column headers in row 1 (D1:W1) so that the
# concatenated "DayHH" strings (e.g. "Mon08") become "Day-H" strings
# (e.g. "Mon-8"), dropping the leading zero on the hour as well.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @{
    "D1" = "Mon-8"
    "E1" = "Mon-10"
    "F1" = "Mon-12"
    "G1" = "Mon-2"
    "H1" = "Tue-8"
    "I1" = "Tue-10"
    "J1" = "Tue-12"
    "K1" = "Tue-2"
    "L1" = "Wed-8"
    "M1" = "Wed-10"
    "N1" = "Wed-12"
    "O1" = "Wed-2"
    "P1" = "Thu-8"
    "Q1" = "Thu-10"
    "R1" = "Thu-12"
    "S1" = "Thu-2"
    "T1" = "Fri-8"
    "U1" = "Fri-10"
    "V1" = "Fri-12"
    "W1" = "Fri-2"
}

foreach ($addr in $headers.Keys) {
    $ws.Range($addr).Value = $headers[$addr]
}

# Update the saved sheet view: scroll so column B is the left-most visible
# column, and select the header row (A1:W1) instead of the single cell I13.
$window = $excel.ActiveWindow
$window.ScrollColumn = 2
$ws.Range("A1:W1").Select()
